$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 557
$ws.Range("F3").Value = 46
$ws.Range("F4").Value = 554
$ws.Range("F6").Value = 1611
$ws.Range("F9").Value = 745
$ws.Range("F10").Value = 2715
$ws.Range("F11").Value = 19
$ws.Range("F12").Value = 1816
$ws.Range("F13").Value = 620
$ws.Range("F14").Value = 300
$ws.Range("F15").Value = 709
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = 6231
$ws.Range("F18").Value = 236
$ws.Range("F21").Value = 3395
$ws.Range("F22").Value = 879
$ws.Range("F24").Value = 74
$ws.Range("F25").Value = 48
$ws.Range("F26").Value = 2463
$ws.Range("F28").Value = 377
$ws.Range("F31").Value = 497
$ws.Range("F32").Value = 1315
$ws.Range("F34").Value = 12
$ws.Range("F35").Value = 82
$ws.Range("F36").Value = 31
$ws.Range("F38").Value = 1489
$ws.Range("F39").Value = 29
$ws.Range("F40").Value = 1443

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 162
$ws.Range("F16").Value = 146
$ws.Range("F17").Value = 337
$ws.Range("F18").Value = 267
$ws.Range("F19").Value = 523

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 922
$ws.Range("F6").Value = 44
$ws.Range("F7").Value = 74

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 557
$ws.Range("F3").Value = 46
$ws.Range("F4").Value = 922
$ws.Range("F6").Value = 554
$ws.Range("F7").Value = 44
$ws.Range("F8").Value = 74
$ws.Range("F16").Value = 2715
$ws.Range("F19").Value = 19
$ws.Range("F21").Value = 1816
$ws.Range("F22").Value = 162
$ws.Range("F23").Value = 620
$ws.Range("F24").Value = 300
$ws.Range("F25").Value = 709
$ws.Range("F26").Value = 6231
$ws.Range("F27").Value = 236
$ws.Range("F30").Value = 3395
$ws.Range("F31").Value = 879
$ws.Range("F33").Value = 74
$ws.Range("F35").Value = 48
$ws.Range("F36").Value = 2464
$ws.Range("F37").Value = 377
$ws.Range("F38").Value = 497
$ws.Range("F39").Value = 1315
$ws.Range("F40").Value = 337
$ws.Range("F41").Value = 267
$ws.Range("F42").Value = 523
$ws.Range("F44").Value = 12
$ws.Range("F45").Value = 82
$ws.Range("F46").Value = 31
$ws.Range("F48").Value = 29
$ws.Range("F50").Value = 1443
